$d = $word.ActiveDocument

$replacements = @(
    @{old = "300÷3=100, 0"; new = "240÷5=48, 0"},
    @{old = "504÷4=126, 0"; new = "362÷3=120, 2"},
    @{old = "458÷5=91, 3"; new = "866÷4=216, 2"},
    @{old = "501÷7=71, 4"; new = "794÷4=198, 2"},
    @{old = "694÷3=231, 1"; new = "325÷5=65, 0"},
    @{old = "731÷2=365, 1"; new = "671÷5=134, 1"},
    @{old = "409÷8=51, 1"; new = "664÷7=94, 6"},
    @{old = "590÷9=65, 5"; new = "848÷3=282, 2"},
    @{old = "695÷2=347, 1"; new = "754÷3=251, 1"},
    @{old = "205÷8=25, 5"; new = "909÷5=181, 4"},
    @{old = "213÷8=26, 5"; new = "110÷2=55, 0"},
    @{old = "130÷7=18, 4"; new = "926÷7=132, 2"},
    @{old = "180÷5=36, 0"; new = "603÷7=86, 1"},
    @{old = "755÷6=125, 5"; new = "622÷8=77, 6"},
    @{old = "763÷9=84, 7"; new = "902÷6=150, 2"},
    @{old = "129÷2=64, 1"; new = "128÷6=21, 2"},
    @{old = "265÷5=53, 0"; new = "790÷5=158, 0"},
    @{old = "932÷3=310, 2"; new = "610÷7=87, 1"},
    @{old = "151÷7=21, 4"; new = "589÷9=65, 4"},
    @{old = "256÷2=128, 0"; new = "829÷7=118, 3"},
    @{old = "518÷8=64, 6"; new = "666÷6=111, 0"},
    @{old = "835÷8=104, 3"; new = "296÷4=74, 0"},
    @{old = "385÷3=128, 1"; new = "756÷7=108, 0"},
    @{old = "782÷7=111, 5"; new = "921÷8=115, 1"},
    @{old = "870÷3=290, 0"; new = "379÷6=63, 1"}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $false, $false, $false, $false,
                         $true, 1, $false, $r.new, 2)
}

Write-Host "Done: applied $($replacements.Count) replacements"
